$wb = $excel.ActiveWorkbook

# Sheet 2: Air_Mystifly_RoundTrip
$ws2 = $wb.Worksheets.Item("Air_Mystifly_RoundTrip")
$ws2.Range("B2").Value = "LOGIN|Search"
$ws2.Activate()
$ws2.Range("B2").Select()

# Sheet 3: Air_Mystifly_Multicity
$ws3 = $wb.Worksheets.Item("Air_Mystifly_Multicity")
$ws3.Range("B2").Value = "LOGIN|Search"
$ws3.Activate()
$ws3.Range("I2").Select()
